$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 7000
$ws.Range("I76").Value = 7000
$ws.Range("J76").Value = 0
$ws.Range("K76").Value = 7000
$ws.Range("L76").Value = 0
$ws.Range("M76").Value = -6685
$ws.Range("N76").ClearContents()

$ws.Range("H79").Value = 7000
$ws.Range("I79").Value = 7000
$ws.Range("J79").Value = 0
$ws.Range("K79").Value = 7000
$ws.Range("L79").Value = 0
$ws.Range("M79").Value = -5908
$ws.Range("N79").ClearContents()

$ws.Range("H132").Value = 6450.5386
$ws.Range("I132").Value = 5344.512
$ws.Range("J132").Value = 10573
$ws.Range("K132").Value = 16033.536
$ws.Range("L132").Value = 31719
$ws.Range("M132").Value = -13503.536
$ws.Range("N132").Value = -36779

$ws.Range("H135").Value = 4746.9
$ws.Range("I135").Value = 2474.1428
$ws.Range("J135").Value = 10050
$ws.Range("K135").Value = 22267.2852
$ws.Range("L135").Value = 90450
$ws.Range("M135").Value = -19732.2852
$ws.Range("N135").Value = -95520

$ws.Range("H138").Value = 2980.5352
$ws.Range("I138").Value = 3783.9333
$ws.Range("J138").Value = 2765.3394
$ws.Range("K138").Value = 11351.7999
$ws.Range("L138").Value = 8296.018199999999
$ws.Range("M138").Value = -6211.7999
$ws.Range("N138").Value = -18576.0182

$ws.Range("H141").Value = 4800
$ws.Range("I141").Value = 4529.35
$ws.Range("K141").Value = 13588.05
$ws.Range("M141").Value = -8408.050000000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4539.7554
$ws.Range("I32").Value = 2657.9863
$ws.Range("J32").Value = 12620.294
$ws.Range("K32").Value = 2657.9863
$ws.Range("L32").Value = 12620.294
$ws.Range("M32").Value = -2370.9863
$ws.Range("N32").Value = -13194.294

$ws.Range("H45").Value = 6917.9165
$ws.Range("I45").Value = 6589
$ws.Range("J45").Value = 7904.6665
$ws.Range("K45").Value = 6589
$ws.Range("L45").Value = 7904.6665
$ws.Range("M45").Value = -6212
$ws.Range("N45").Value = -8658.666499999999

$ws.Range("H102").Value = 11021.19
$ws.Range("I102").Value = 1458.0714
$ws.Range("J102").Value = 30147.428
$ws.Range("K102").Value = 1458.0714
$ws.Range("L102").Value = 30147.428
$ws.Range("M102").Value = 163.9286
$ws.Range("N102").Value = -33391.428

$ws.Range("H122").Value = 2596603.2
$ws.Range("I122").Value = 13820652
$ws.Range("J122").Value = 6438.231
$ws.Range("K122").Value = 41461956
$ws.Range("L122").Value = 19314.693
$ws.Range("M122").Value = -41459506
$ws.Range("N122").Value = -24214.693

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 1524.6897
$ws.Range("I107").Value = 1588.2084
$ws.Range("J107").Value = 1219.8
$ws.Range("K107").Value = 1588.2084
$ws.Range("L107").Value = 1219.8
$ws.Range("M107").Value = 331.7916
$ws.Range("N107").Value = -5059.8

$ws.Range("H128").Value = 7848.6
$ws.Range("I128").Value = 7848.6
$ws.Range("K128").Value = 23545.8
$ws.Range("M128").Value = -21055.8

$ws.Range("H132").Value = 82363.86
$ws.Range("J132").Value = 82363.86
$ws.Range("L132").Value = 82363.86
$ws.Range("N132").Value = -92483.86

$ws.Range("H134").Value = 26077.5
$ws.Range("I134").Value = 24085.148
$ws.Range("J134").Value = 44805.6
$ws.Range("K134").Value = 72255.444
$ws.Range("L134").Value = 134416.8
$ws.Range("M134").Value = -69720.444
$ws.Range("N134").Value = -139486.8

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 10149.533
$ws.Range("I16").Value = 8827.875
$ws.Range("J16").Value = 11660
$ws.Range("K16").Value = 8827.875
$ws.Range("L16").Value = 11660
$ws.Range("M16").Value = -8540.875
$ws.Range("N16").Value = -12234

$ws.Range("H22").Value = 1154.375
$ws.Range("I22").Value = 782.3333
$ws.Range("J22").Value = 1774.4445
$ws.Range("K22").Value = 782.3333
$ws.Range("L22").Value = 1774.4445
$ws.Range("M22").Value = -432.3333
$ws.Range("N22").Value = -2474.4445

$ws.Range("H31").Value = 56426.047
$ws.Range("I31").Value = 88974.836
$ws.Range("J31").Value = 17367.5
$ws.Range("K31").Value = 88974.836
$ws.Range("L31").Value = 17367.5
$ws.Range("M31").Value = -88679.836
$ws.Range("N31").Value = -17957.5

$ws.Range("H34").Value = 56426.047
$ws.Range("I34").Value = 88974.836
$ws.Range("J34").Value = 17367.5
$ws.Range("K34").Value = 88974.836
$ws.Range("L34").Value = 17367.5
$ws.Range("M34").Value = -88772.836
$ws.Range("N34").Value = -17771.5

$ws.Range("H113").Value = 10149.533
$ws.Range("I113").Value = 8827.875
$ws.Range("J113").Value = 11660
$ws.Range("K113").Value = 8827.875
$ws.Range("L113").Value = 11660
$ws.Range("M113").Value = -6657.875
$ws.Range("N113").Value = -16000

$ws.Range("H132").Value = 3345.1355
$ws.Range("I132").Value = 1178.5094
$ws.Range("J132").Value = 22483.666
$ws.Range("K132").Value = 3535.5282
$ws.Range("L132").Value = 67450.99800000001
$ws.Range("M132").Value = -1005.5282
$ws.Range("N132").Value = -72510.99800000001

$ws.Range("H134").Value = 3611.106
$ws.Range("I134").Value = 1437.1777
$ws.Range("J134").Value = 8269.522999999999
$ws.Range("K134").Value = 4311.5331
$ws.Range("L134").Value = 24808.569
$ws.Range("M134").Value = -1776.5331
$ws.Range("N134").Value = -29878.569

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 1119.6428
$ws.Range("I68").Value = 841.1429000000001
$ws.Range("K68").Value = 2523.4287
$ws.Range("M68").Value = -1712.4287

$ws.Range("H71").Value = 1119.6428
$ws.Range("I71").Value = 841.1429000000001
$ws.Range("K71").Value = 7570.2861
$ws.Range("M71").Value = -3514.2861

$ws.Range("H86").Value = 702.5417
$ws.Range("I86").Value = 749.1053000000001
$ws.Range("J86").Value = 525.6
$ws.Range("K86").Value = 2247.3159
$ws.Range("L86").Value = 1576.8
$ws.Range("M86").Value = -1061.3159
$ws.Range("N86").Value = -3948.8

$ws.Range("H89").Value = 702.5417
$ws.Range("I89").Value = 749.1053000000001
$ws.Range("J89").Value = 525.6
$ws.Range("K89").Value = 6741.947700000001
$ws.Range("L89").Value = 4730.400000000001
$ws.Range("M89").Value = -813.9477000000006
$ws.Range("N89").Value = -16586.4

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 30598.125
$ws.Range("I80").Value = 33962.668
$ws.Range("K80").Value = 33962.668
$ws.Range("M80").Value = -32964.668

$ws.Range("H83").Value = 30598.125
$ws.Range("I83").Value = 33962.668
$ws.Range("K83").Value = 169813.34
$ws.Range("M83").Value = -164821.34

$ws.Range("H102").Value = 9656992
$ws.Range("I102").Value = 9656992
$ws.Range("K102").Value = 9656992
$ws.Range("M102").Value = -9655370

$ws.Range("H113").Value = 5436.0713
$ws.Range("I113").Value = 5506.2856
$ws.Range("J113").Value = 5365.857
$ws.Range("K113").Value = 5506.2856
$ws.Range("L113").Value = 5365.857
$ws.Range("M113").Value = -3336.2856
$ws.Range("N113").Value = -9705.857

$ws.Range("H123").Value = 56449.25
$ws.Range("J123").Value = 56449.25
$ws.Range("L123").Value = 56449.25
$ws.Range("N123").Value = -61349.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2988034
$ws.Range("I7").Value = 5967319
$ws.Range("K7").Value = 5967319
$ws.Range("M7").Value = -5967207

$ws.Range("H40").Value = 2807607.2
$ws.Range("I40").Value = 4959.75
$ws.Range("J40").Value = 6544471
$ws.Range("K40").Value = 4959.75
$ws.Range("L40").Value = 6544471
$ws.Range("M40").Value = -4823.75
$ws.Range("N40").Value = -6544743

$ws.Range("H93").Value = 50010844
$ws.Range("I93").Value = 100011490
$ws.Range("K93").Value = 100011490
$ws.Range("M93").Value = -100010242

$ws.Range("H126").Value = 2988034
$ws.Range("I126").Value = 5967319
$ws.Range("K126").Value = 17901957
$ws.Range("M126").Value = -17899487

$ws.Range("H132").Value = 6338.357
$ws.Range("I132").Value = 2592.2666
$ws.Range("K132").Value = 7776.7998
$ws.Range("M132").Value = -5246.7998
